$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 19: KODEX 미국달러선물레버리지 / 2018-10-11-목 / short / 9190 / 110 / 2%
$ws.Cells.Item(19,1).Value = "KODEX 미국달러선물레버리지"
$ws.Cells.Item(19,2).Value = "2018-10-11-목"
$ws.Cells.Item(19,3).Value = "short"
$ws.Cells.Item(19,4).Value = 9190
$ws.Cells.Item(19,4).NumberFormat = "#,##0"
$ws.Cells.Item(19,5).Value = 110
$ws.Cells.Item(19,6).Value = 0.02
$ws.Cells.Item(19,6).NumberFormat = "0%"

# New row 20: KODEX 미국달러선물레버리지 / 2018-10-11-목 / short / 9200 / 110 / 3%
$ws.Cells.Item(20,1).Value = "KODEX 미국달러선물레버리지"
$ws.Cells.Item(20,2).Value = "2018-10-11-목"
$ws.Cells.Item(20,3).Value = "short"
$ws.Cells.Item(20,4).Value = 9200
$ws.Cells.Item(20,4).NumberFormat = "#,##0"
$ws.Cells.Item(20,5).Value = 110
$ws.Cells.Item(20,6).Value = 0.03
$ws.Cells.Item(20,6).NumberFormat = "0%"

# Update selected cell shown in the sheet view
$ws.Range("F21").Select()
